$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two leading blank columns (old A, old B) so old C:T shifts to new A:R.
$ws.Columns("A:B").Delete()

# The first column (new A, old C) now holds "1. 会社名" and needs a wider column (17 chars).
$ws.Columns("A:A").ColumnWidth = 16.285714285714285

# Split "1. 会社名" into two runs: "1. " (default font) and "会社名" (MS Pゴシック).
$ws.Range("A1").Value = "1. 会社名"
$chars = $ws.Range("A1").Characters(4, 3)
$chars.Font.Name = "ＭＳ Ｐゴシック"
$chars.Font.Size = 11
$chars.Font.Charset = 128

# Update the print area to match the new (narrower) layout.
$ws.PageSetup.PrintArea = '$A$1:$R$31'

Write-Output "done"
